# Update "想去人数" (F column) figures across the workbook's sheets to
# reflect newly generated output, as published to gh-pages.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6609
$ws1.Range("F3").Value  = 776
$ws1.Range("F5").Value  = 121
$ws1.Range("F6").Value  = 651
$ws1.Range("F7").Value  = 215
$ws1.Range("F8").Value  = 47
$ws1.Range("F10").Value = 1275
$ws1.Range("F12").Value = 100
$ws1.Range("F13").Value = 520
$ws1.Range("F15").Value = 359
$ws1.Range("F16").Value = 1045
$ws1.Range("F17").Value = 1460
$ws1.Range("F24").Value = 2304

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 85

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1610

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1610
$ws4.Range("F8").Value  = 6609
$ws4.Range("F10").Value = 776
$ws4.Range("F12").Value = 121
$ws4.Range("F13").Value = 651
$ws4.Range("F14").Value = 651
$ws4.Range("F15").Value = 215
$ws4.Range("F16").Value = 47
$ws4.Range("F23").Value = 85
$ws4.Range("F24").Value = 1275
$ws4.Range("F26").Value = 100
$ws4.Range("F27").Value = 520
$ws4.Range("F32").Value = 359
$ws4.Range("F33").Value = 1045
$ws4.Range("F34").Value = 1460
$ws4.Range("F43").Value = 2304
